$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attendance")

$attendance = @{
    2  = 22
    3  = 10
    4  = 15
    5  = 4
    6  = 36
    7  = 12
    8  = 22
    9  = 8
    10 = 15
    11 = 20
    12 = 10
    13 = 25
    14 = 10
    15 = 23
    16 = 15
    17 = 10
    18 = 10
    19 = 12
    20 = 17
    21 = 10
}

foreach ($row in ($attendance.Keys | Sort-Object)) {
    $ws.Cells.Item($row, 4).Value = $attendance[$row]
}

$ws.Range("D18").Select() | Out-Null
